# Capstone Final Project Presentation - "final final final version" edit
#
# 1) Swap the order of slides 10 and 11 (the "Potentially Faulty
#    Sensors/Unexpected Patters in Sensor Behavior" slide and the
#    "Significance Cut-Off Calculations" slide) so the Significance
#    Cut-Off slide now comes right after slide 9, and the Potentially
#    Faulty Sensors slide follows it.
# 2) Tighten up the wording of the three guiding questions on slide 4
#    (change "you" -> "we", drop the "Now turn your attention..." lead
#    in, and simplify the "factories" question) and let PowerPoint
#    recompute the shrink-to-fit font scale for that text box.

$p = $ppt.ActivePresentation

# --- 1) Reorder slides 10 and 11 ------------------------------------
$s11 = $p.Slides.Item(11)
$s11.MoveTo(10)

# --- 2) Slide 4 question wording -------------------------------------
$s4 = $p.Slides.Item(4)
$body = $s4.Shapes.Item(5)

$body.TextFrame.TextRange.Paragraphs(1).Text = "Characterize the sensors" + [char]0x2019 + " performance and operation. Are they all working properly at all times? Can we detect any unexpected behaviors of the sensors through analyzing the readings they capture?"

$body.TextFrame.TextRange.Paragraphs(3).Text = "Which chemicals are being detected by the sensor group? What patterns of chemical releases do we see, as being reported in the data?"

$body.TextFrame.TextRange.Paragraphs(5).Text = "Which factories are responsible for which chemical releases? For the identified factories, are there any observed patterns of operation that can be described from the data?"

# Let PowerPoint recompute the shrink-to-fit scale for the new text
# (drops the stale fontScale="92500" left over from the old wording).
$body.TextFrame.AutoSize = 2
